$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "negative power is now zero" - set min Leistung (column H) rows 3-6 to 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
